# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 0.1190320826869504; C = 1.655778082260271;  D = 0.1494219747398047; E = 10.19245300693656 }
    3 = @{ B = 1.455362044514542;  C = 10.34677158129881;  D = 0.7527432677738641; E = 10.19245300693656 }
    4 = @{ B = 3.286832544864788;  C = 1.655778082260271;  D = 3.537761648806719;  E = 0.4942365360607697 }
    5 = @{ B = 3.286832544864788;  C = 1.655778082260271;  D = 0.1494219747398047; E = 0.4942365360607697 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    # F column (Win) is unchanged; G (sum) = B + C + D + E
    $sum = $vals.B + $vals.C + $vals.D + $vals.E
    $ws.Range("G$row").Value = $sum
}
